$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01293466051926884
$ws.Range("C2").Value = 1.689667739057654 * [Math]::Pow(10, -9)
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1136.587612837883
